$wb = $excel.ActiveWorkbook

# --- Summary sheet: Total Trades 12 -> 13, Win Rate % 33.33 -> 30.77 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 13
$wsSummary.Range("B9").Value = 30.77

# --- Strategy Status sheet: MarketMaking Trades 12 -> 13, Win Rate % 33.33 -> 30.77 ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 13
$wsStatus.Range("G4").Value = 30.77

# --- New closed trade (#13) appended to "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A14").Value = 13

    # Date/time columns are stored as plain text in this workbook (not Excel
    # date serials). Force text formatting before assigning so the "YYYY-MM-DD"
    # string isn't auto-converted into a date serial, then clear the formatting
    # again so the cell keeps the original (unstyled) look - only the stored
    # value/type is what we actually need to change.
    $ws.Range("B14").NumberFormat = "@"
    $ws.Range("B14").Value = "2026-02-17"
    $ws.Range("B14").ClearFormats()

    $ws.Range("C14").Value = "04:07:19"

    $ws.Range("D14").Value = "MarketMaking"
    $ws.Range("E14").Value = "DOWN"
    $ws.Range("F14").Value = 0.8
    $ws.Range("G14").Value = 0.8
    $ws.Range("H14").Value = "CLOSED"
    $ws.Range("I14").Value = 0
    $ws.Range("J14").Value = 0
    $ws.Range("K14").Value = 100.02
    $ws.Range("L14").Value = 0
    $ws.Range("M14").Value = 0
    $ws.Range("N14").Value = 0.6
    $ws.Range("O14").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P14").Value = "early_exit"
    $ws.Range("Q14").Value = 0.12
}
